# Added Audio Amp shutdown function.
#
# Semantic changes applied to the "ToDos" sheet (Sheet1):
#   - Row 31 "Low Power mode": Status Open -> Ongoing, Notes set to the new
#     "Audio Amp shuts down when not in use" item.
#   - Rows 42-44: Status "Done" -> "Closed" (the "Done" status value is
#     retired / no longer used anywhere in the workbook).
#   - Row 47 "LEDs keep flashing after reset": Status Open -> Closed.
#   - Row 78 "Hissing, crackling, high pitch noise come from speakers in rest":
#     Status Open -> Ongoing, Notes set to the new Audio Amp item (same new
#     todo note as row 31).
#   - Row 79 "Intercom2 speaker does not work": Status Open -> Closed.
#   - The sheet view scrolled down / selection moved to C80 (bottom of the
#     list), reflecting that the edits were made near the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$note = "Audio Amp shuts down when not in use"

# Row 31 - Low Power mode
$ws.Range("C31").Value = "Ongoing"
$ws.Range("D31").Value = $note

# Rows 42-44 - Done -> Closed
$ws.Range("C42").Value = "Closed"
$ws.Range("C43").Value = "Closed"
$ws.Range("C44").Value = "Closed"

# Row 47 - LEDs keep flashing after reset
$ws.Range("C47").Value = "Closed"

# Row 78 - Hissing, crackling, high pitch noise come from speakers in rest
$ws.Range("C78").Value = "Ongoing"
$ws.Range("D78").Value = $note

# Row 79 - Intercom2 speaker does not work
$ws.Range("C79").Value = "Closed"

# Match the author's final view state: scrolled to bottom of the list with
# C80 selected.
$ws.Range("C80").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 52
